$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("snapshot")
$ws2 = $wb.Worksheets.Item("returned")
$ws3 = $wb.Worksheets.Item("new_injured")

# --- Sheet1 "snapshot": rewrite data rows 2..33 (row 34 removed) ---
$ws1.Cells.Item(2, 1).Value = 'АВТ'
$ws1.Cells.Item(2, 2).Value = 'Автомобилист'
$ws1.Cells.Item(2, 3).Value = 'avtomobilist'
$ws1.Cells.Item(2, 4).Value = 'Зборовский Сергей'
$ws1.Cells.Item(2, 5).Value = '2'
$ws1.Cells.Item(2, 6).Value = 'защитник'
$ws1.Cells.Item(2, 7).Value = '20989'
$ws1.Cells.Item(2, 8).Value = '1369_АВТ_зборовскийсергей'
$ws1.Cells.Item(2, 9).Value = 'injured_active'
$ws1.Cells.Item(2, 10).Value = 'https://www.khl.ru/clubs/avtomobilist/team/'
$ws1.Cells.Item(2, 11).Value = '2025-11-29T03:01:44.128873+00:00'

$ws1.Cells.Item(3, 1).Value = 'АДМ'
$ws1.Cells.Item(3, 2).Value = 'Адмирал'
$ws1.Cells.Item(3, 3).Value = 'admiral'
$ws1.Cells.Item(3, 4).Value = 'Грман Марио'
$ws1.Cells.Item(3, 5).Value = '77'
$ws1.Cells.Item(3, 6).Value = 'защитник'
$ws1.Cells.Item(3, 7).Value = '31232'
$ws1.Cells.Item(3, 8).Value = '1369_АДМ_грманмарио'
$ws1.Cells.Item(3, 9).Value = 'injured_active'
$ws1.Cells.Item(3, 10).Value = 'https://www.khl.ru/clubs/admiral/team/'
$ws1.Cells.Item(3, 11).Value = '2025-11-29T03:01:46.486325+00:00'

$ws1.Cells.Item(4, 1).Value = 'АДМ'
$ws1.Cells.Item(4, 2).Value = 'Адмирал'
$ws1.Cells.Item(4, 3).Value = 'admiral'
$ws1.Cells.Item(4, 4).Value = 'Дарьин Александр'
$ws1.Cells.Item(4, 5).Value = '41'
$ws1.Cells.Item(4, 6).Value = 'нападающий'
$ws1.Cells.Item(4, 7).Value = '26282'
$ws1.Cells.Item(4, 8).Value = '1369_АДМ_дарьиналександр'
$ws1.Cells.Item(4, 9).Value = 'injured_active'
$ws1.Cells.Item(4, 10).Value = 'https://www.khl.ru/clubs/admiral/team/'
$ws1.Cells.Item(4, 11).Value = '2025-11-29T03:01:46.486357+00:00'

$ws1.Cells.Item(5, 1).Value = 'АДМ'
$ws1.Cells.Item(5, 2).Value = 'Адмирал'
$ws1.Cells.Item(5, 3).Value = 'admiral'
$ws1.Cells.Item(5, 4).Value = 'Сошников Никита'
$ws1.Cells.Item(5, 5).Value = '90'
$ws1.Cells.Item(5, 6).Value = 'нападающий'
$ws1.Cells.Item(5, 7).Value = '16731'
$ws1.Cells.Item(5, 8).Value = '1369_АДМ_сошниковникита'
$ws1.Cells.Item(5, 9).Value = 'injured_active'
$ws1.Cells.Item(5, 10).Value = 'https://www.khl.ru/clubs/admiral/team/'
$ws1.Cells.Item(5, 11).Value = '2025-11-29T03:01:46.486404+00:00'

$ws1.Cells.Item(6, 1).Value = 'АКБ'
$ws1.Cells.Item(6, 2).Value = 'Ак Барс'
$ws1.Cells.Item(6, 3).Value = 'ak_bars'
$ws1.Cells.Item(6, 4).Value = 'Яруллин Альберт'
$ws1.Cells.Item(6, 5).Value = '33'
$ws1.Cells.Item(6, 6).Value = 'защитник'
$ws1.Cells.Item(6, 7).Value = '16365'
$ws1.Cells.Item(6, 8).Value = '1369_АКБ_яруллинальберт'
$ws1.Cells.Item(6, 9).Value = 'injured_active'
$ws1.Cells.Item(6, 10).Value = 'https://www.khl.ru/clubs/ak_bars/team/'
$ws1.Cells.Item(6, 11).Value = '2025-11-29T03:01:49.291067+00:00'

$ws1.Cells.Item(7, 1).Value = 'АМР'
$ws1.Cells.Item(7, 2).Value = 'Амур'
$ws1.Cells.Item(7, 3).Value = 'amur'
$ws1.Cells.Item(7, 4).Value = 'Абросимов Роман'
$ws1.Cells.Item(7, 5).Value = '94'
$ws1.Cells.Item(7, 6).Value = 'защитник'
$ws1.Cells.Item(7, 7).Value = '17968'
$ws1.Cells.Item(7, 8).Value = '1369_АМР_абросимовроман'
$ws1.Cells.Item(7, 9).Value = 'injured_active'
$ws1.Cells.Item(7, 10).Value = 'https://www.khl.ru/clubs/amur/team/'
$ws1.Cells.Item(7, 11).Value = '2025-11-29T03:01:52.246627+00:00'

$ws1.Cells.Item(8, 1).Value = 'БАР'
$ws1.Cells.Item(8, 2).Value = 'Барыс'
$ws1.Cells.Item(8, 3).Value = 'barys'
$ws1.Cells.Item(8, 4).Value = 'Уотерспун Тайлер'
$ws1.Cells.Item(8, 5).Value = '26'
$ws1.Cells.Item(8, 6).Value = 'защитник'
$ws1.Cells.Item(8, 7).Value = '45769'
$ws1.Cells.Item(8, 8).Value = '1369_БАР_уотерспунтайлер'
$ws1.Cells.Item(8, 9).Value = 'injured_active'
$ws1.Cells.Item(8, 10).Value = 'https://www.khl.ru/clubs/barys/team/'
$ws1.Cells.Item(8, 11).Value = '2025-11-29T03:01:55.054289+00:00'

$ws1.Cells.Item(9, 1).Value = 'ДИН'
$ws1.Cells.Item(9, 2).Value = 'Динамо М'
$ws1.Cells.Item(9, 3).Value = 'dynamo_msk'
$ws1.Cells.Item(9, 4).Value = 'Готовец Кирилл'
$ws1.Cells.Item(9, 5).Value = '41'
$ws1.Cells.Item(9, 6).Value = 'защитник'
$ws1.Cells.Item(9, 7).Value = '16034'
$ws1.Cells.Item(9, 8).Value = '1369_ДИН_готовецкирилл'
$ws1.Cells.Item(9, 9).Value = 'injured_active'
$ws1.Cells.Item(9, 10).Value = 'https://www.khl.ru/clubs/dynamo_msk/team/'
$ws1.Cells.Item(9, 11).Value = '2025-11-29T03:01:57.858109+00:00'

$ws1.Cells.Item(10, 1).Value = 'ЛАД'
$ws1.Cells.Item(10, 2).Value = 'Лада'
$ws1.Cells.Item(10, 3).Value = 'lada'
$ws1.Cells.Item(10, 4).Value = 'Обидин Андрей'
$ws1.Cells.Item(10, 5).Value = '58'
$ws1.Cells.Item(10, 6).Value = 'нападающий'
$ws1.Cells.Item(10, 7).Value = '21310'
$ws1.Cells.Item(10, 8).Value = '1369_ЛАД_обидинандрей'
$ws1.Cells.Item(10, 9).Value = 'injured_active'
$ws1.Cells.Item(10, 10).Value = 'https://www.khl.ru/clubs/lada/team/'
$ws1.Cells.Item(10, 11).Value = '2025-11-29T03:02:02.904222+00:00'

$ws1.Cells.Item(11, 1).Value = 'ЛАД'
$ws1.Cells.Item(11, 2).Value = 'Лада'
$ws1.Cells.Item(11, 3).Value = 'lada'
$ws1.Cells.Item(11, 4).Value = 'Ожгихин Алексей'
$ws1.Cells.Item(11, 5).Value = '43'
$ws1.Cells.Item(11, 6).Value = 'нападающий'
$ws1.Cells.Item(11, 7).Value = '23021'
$ws1.Cells.Item(11, 8).Value = '1369_ЛАД_ожгихиналексей'
$ws1.Cells.Item(11, 9).Value = 'injured_active'
$ws1.Cells.Item(11, 10).Value = 'https://www.khl.ru/clubs/lada/team/'
$ws1.Cells.Item(11, 11).Value = '2025-11-29T03:02:02.904252+00:00'

$ws1.Cells.Item(12, 1).Value = 'ЛОК'
$ws1.Cells.Item(12, 2).Value = 'Локомотив'
$ws1.Cells.Item(12, 3).Value = 'lokomotiv'
$ws1.Cells.Item(12, 4).Value = 'Сергеев Андрей'
$ws1.Cells.Item(12, 5).Value = '99'
$ws1.Cells.Item(12, 6).Value = 'защитник'
$ws1.Cells.Item(12, 7).Value = '15416'
$ws1.Cells.Item(12, 8).Value = '1369_ЛОК_сергеевандрей'
$ws1.Cells.Item(12, 9).Value = 'injured_active'
$ws1.Cells.Item(12, 10).Value = 'https://www.khl.ru/clubs/lokomotiv/team/'
$ws1.Cells.Item(12, 11).Value = '2025-11-29T03:02:05.740750+00:00'

$ws1.Cells.Item(13, 1).Value = 'ММГ'
$ws1.Cells.Item(13, 2).Value = 'Металлург Мг'
$ws1.Cells.Item(13, 3).Value = 'metallurg_mg'
$ws1.Cells.Item(13, 4).Value = 'Козлов Андрей Е'
$ws1.Cells.Item(13, 5).Value = '39'
$ws1.Cells.Item(13, 6).Value = 'нападающий'
$ws1.Cells.Item(13, 7).Value = '40899'
$ws1.Cells.Item(13, 8).Value = '1369_ММГ_козловандрейе'
$ws1.Cells.Item(13, 9).Value = 'injured_active'
$ws1.Cells.Item(13, 10).Value = 'https://www.khl.ru/clubs/metallurg_mg/team/'
$ws1.Cells.Item(13, 11).Value = '2025-11-29T03:02:08.439163+00:00'

$ws1.Cells.Item(14, 1).Value = 'ММГ'
$ws1.Cells.Item(14, 2).Value = 'Металлург Мг'
$ws1.Cells.Item(14, 3).Value = 'metallurg_mg'
$ws1.Cells.Item(14, 4).Value = 'Сиряцкий Александр'
$ws1.Cells.Item(14, 5).Value = '74'
$ws1.Cells.Item(14, 6).Value = 'защитник'
$ws1.Cells.Item(14, 7).Value = '42458'
$ws1.Cells.Item(14, 8).Value = '1369_ММГ_сиряцкийалександр'
$ws1.Cells.Item(14, 9).Value = 'injured_active'
$ws1.Cells.Item(14, 10).Value = 'https://www.khl.ru/clubs/metallurg_mg/team/'
$ws1.Cells.Item(14, 11).Value = '2025-11-29T03:02:08.439194+00:00'

$ws1.Cells.Item(15, 1).Value = 'НХК'
$ws1.Cells.Item(15, 2).Value = 'Нефтехимик'
$ws1.Cells.Item(15, 3).Value = 'neftekhimik'
$ws1.Cells.Item(15, 4).Value = 'Хлыстов Никита'
$ws1.Cells.Item(15, 5).Value = '7'
$ws1.Cells.Item(15, 6).Value = 'защитник'
$ws1.Cells.Item(15, 7).Value = '17881'
$ws1.Cells.Item(15, 8).Value = '1369_НХК_хлыстовникита'
$ws1.Cells.Item(15, 9).Value = 'injured_active'
$ws1.Cells.Item(15, 10).Value = 'https://www.khl.ru/clubs/neftekhimik/team/'
$ws1.Cells.Item(15, 11).Value = '2025-11-29T03:02:11.204790+00:00'

$ws1.Cells.Item(16, 1).Value = 'СЕВ'
$ws1.Cells.Item(16, 2).Value = 'Северсталь'
$ws1.Cells.Item(16, 3).Value = 'severstal'
$ws1.Cells.Item(16, 4).Value = 'Квочко Илья'
$ws1.Cells.Item(16, 5).Value = '51'
$ws1.Cells.Item(16, 6).Value = 'нападающий'
$ws1.Cells.Item(16, 7).Value = '39445'
$ws1.Cells.Item(16, 8).Value = '1369_СЕВ_квочкоилья'
$ws1.Cells.Item(16, 9).Value = 'injured_active'
$ws1.Cells.Item(16, 10).Value = 'https://www.khl.ru/clubs/severstal/team/'
$ws1.Cells.Item(16, 11).Value = '2025-11-29T03:02:14.048242+00:00'

$ws1.Cells.Item(17, 1).Value = 'СЕВ'
$ws1.Cells.Item(17, 2).Value = 'Северсталь'
$ws1.Cells.Item(17, 3).Value = 'severstal'
$ws1.Cells.Item(17, 4).Value = 'Фомин Макар'
$ws1.Cells.Item(17, 5).Value = '77'
$ws1.Cells.Item(17, 6).Value = 'защитник'
$ws1.Cells.Item(17, 7).Value = '42087'
$ws1.Cells.Item(17, 8).Value = '1369_СЕВ_фоминмакар'
$ws1.Cells.Item(17, 9).Value = 'injured_active'
$ws1.Cells.Item(17, 10).Value = 'https://www.khl.ru/clubs/severstal/team/'
$ws1.Cells.Item(17, 11).Value = '2025-11-29T03:02:14.048274+00:00'

$ws1.Cells.Item(18, 1).Value = 'СЕВ'
$ws1.Cells.Item(18, 2).Value = 'Северсталь'
$ws1.Cells.Item(18, 3).Value = 'severstal'
$ws1.Cells.Item(18, 4).Value = 'Цицюра Владислав'
$ws1.Cells.Item(18, 5).Value = '10'
$ws1.Cells.Item(18, 6).Value = 'нападающий'
$ws1.Cells.Item(18, 7).Value = '23840'
$ws1.Cells.Item(18, 8).Value = '1369_СЕВ_цицюравладислав'
$ws1.Cells.Item(18, 9).Value = 'injured_active'
$ws1.Cells.Item(18, 10).Value = 'https://www.khl.ru/clubs/severstal/team/'
$ws1.Cells.Item(18, 11).Value = '2025-11-29T03:02:14.048293+00:00'

$ws1.Cells.Item(19, 1).Value = 'СИБ'
$ws1.Cells.Item(19, 2).Value = 'Сибирь'
$ws1.Cells.Item(19, 3).Value = 'sibir'
$ws1.Cells.Item(19, 4).Value = 'Аланов Егор'
$ws1.Cells.Item(19, 5).Value = '21'
$ws1.Cells.Item(19, 6).Value = 'защитник'
$ws1.Cells.Item(19, 7).Value = '26698'
$ws1.Cells.Item(19, 8).Value = '1369_СИБ_алановегор'
$ws1.Cells.Item(19, 9).Value = 'injured_active'
$ws1.Cells.Item(19, 10).Value = 'https://www.khl.ru/clubs/sibir/team/'
$ws1.Cells.Item(19, 11).Value = '2025-11-29T03:02:16.424683+00:00'

$ws1.Cells.Item(20, 1).Value = 'СИБ'
$ws1.Cells.Item(20, 2).Value = 'Сибирь'
$ws1.Cells.Item(20, 3).Value = 'sibir'
$ws1.Cells.Item(20, 4).Value = 'Пьянов Валентин'
$ws1.Cells.Item(20, 5).Value = '45'
$ws1.Cells.Item(20, 6).Value = 'нападающий'
$ws1.Cells.Item(20, 7).Value = '16195'
$ws1.Cells.Item(20, 8).Value = '1369_СИБ_пьяноввалентин'
$ws1.Cells.Item(20, 9).Value = 'injured_active'
$ws1.Cells.Item(20, 10).Value = 'https://www.khl.ru/clubs/sibir/team/'
$ws1.Cells.Item(20, 11).Value = '2025-11-29T03:02:16.424718+00:00'

$ws1.Cells.Item(21, 1).Value = 'СИБ'
$ws1.Cells.Item(21, 2).Value = 'Сибирь'
$ws1.Cells.Item(21, 3).Value = 'sibir'
$ws1.Cells.Item(21, 4).Value = 'Широков Сергей'
$ws1.Cells.Item(21, 5).Value = '52'
$ws1.Cells.Item(21, 6).Value = 'нападающий'
$ws1.Cells.Item(21, 7).Value = '524'
$ws1.Cells.Item(21, 8).Value = '1369_СИБ_широковсергей'
$ws1.Cells.Item(21, 9).Value = 'injured_active'
$ws1.Cells.Item(21, 10).Value = 'https://www.khl.ru/clubs/sibir/team/'
$ws1.Cells.Item(21, 11).Value = '2025-11-29T03:02:16.424737+00:00'

$ws1.Cells.Item(22, 1).Value = 'СКА'
$ws1.Cells.Item(22, 2).Value = 'СКА'
$ws1.Cells.Item(22, 3).Value = 'ska'
$ws1.Cells.Item(22, 4).Value = 'Зайцев Никита И'
$ws1.Cells.Item(22, 5).Value = '22'
$ws1.Cells.Item(22, 6).Value = 'защитник'
$ws1.Cells.Item(22, 7).Value = '16024'
$ws1.Cells.Item(22, 8).Value = '1369_СКА_зайцевникитаи'
$ws1.Cells.Item(22, 9).Value = 'injured_active'
$ws1.Cells.Item(22, 10).Value = 'https://www.khl.ru/clubs/ska/team/'
$ws1.Cells.Item(22, 11).Value = '2025-11-29T03:02:18.715622+00:00'

$ws1.Cells.Item(23, 1).Value = 'СКА'
$ws1.Cells.Item(23, 2).Value = 'СКА'
$ws1.Cells.Item(23, 3).Value = 'ska'
$ws1.Cells.Item(23, 4).Value = 'Короткий Матвей'
$ws1.Cells.Item(23, 5).Value = '71'
$ws1.Cells.Item(23, 6).Value = 'нападающий'
$ws1.Cells.Item(23, 7).Value = '41566'
$ws1.Cells.Item(23, 8).Value = '1369_СКА_короткийматвей'
$ws1.Cells.Item(23, 9).Value = 'injured_active'
$ws1.Cells.Item(23, 10).Value = 'https://www.khl.ru/clubs/ska/team/'
$ws1.Cells.Item(23, 11).Value = '2025-11-29T03:02:18.715654+00:00'

$ws1.Cells.Item(24, 1).Value = 'СОЧ'
$ws1.Cells.Item(24, 2).Value = 'ХК Сочи'
$ws1.Cells.Item(24, 3).Value = 'hc_sochi'
$ws1.Cells.Item(24, 4).Value = 'Бикмуллин Рафаэль'
$ws1.Cells.Item(24, 5).Value = '24'
$ws1.Cells.Item(24, 6).Value = 'нападающий'
$ws1.Cells.Item(24, 7).Value = '22424'
$ws1.Cells.Item(24, 8).Value = '1369_СОЧ_бикмуллинрафаэль'
$ws1.Cells.Item(24, 9).Value = 'injured_active'
$ws1.Cells.Item(24, 10).Value = 'https://www.khl.ru/clubs/hc_sochi/team/'
$ws1.Cells.Item(24, 11).Value = '2025-11-29T03:02:21.208173+00:00'

$ws1.Cells.Item(25, 1).Value = 'СОЧ'
$ws1.Cells.Item(25, 2).Value = 'ХК Сочи'
$ws1.Cells.Item(25, 3).Value = 'hc_sochi'
$ws1.Cells.Item(25, 4).Value = 'Венгрыжановский Денис'
$ws1.Cells.Item(25, 5).Value = '9'
$ws1.Cells.Item(25, 6).Value = 'нападающий'
$ws1.Cells.Item(25, 7).Value = '31892'
$ws1.Cells.Item(25, 8).Value = '1369_СОЧ_венгрыжановскийденис'
$ws1.Cells.Item(25, 9).Value = 'injured_active'
$ws1.Cells.Item(25, 10).Value = 'https://www.khl.ru/clubs/hc_sochi/team/'
$ws1.Cells.Item(25, 11).Value = '2025-11-29T03:02:21.208204+00:00'

$ws1.Cells.Item(26, 1).Value = 'СОЧ'
$ws1.Cells.Item(26, 2).Value = 'ХК Сочи'
$ws1.Cells.Item(26, 3).Value = 'hc_sochi'
$ws1.Cells.Item(26, 4).Value = 'Сушко Илья'
$ws1.Cells.Item(26, 5).Value = '51'
$ws1.Cells.Item(26, 6).Value = 'защитник'
$ws1.Cells.Item(26, 7).Value = '22691'
$ws1.Cells.Item(26, 8).Value = '1369_СОЧ_сушкоилья'
$ws1.Cells.Item(26, 9).Value = 'injured_active'
$ws1.Cells.Item(26, 10).Value = 'https://www.khl.ru/clubs/hc_sochi/team/'
$ws1.Cells.Item(26, 11).Value = '2025-11-29T03:02:21.208220+00:00'

$ws1.Cells.Item(27, 1).Value = 'СОЧ'
$ws1.Cells.Item(27, 2).Value = 'ХК Сочи'
$ws1.Cells.Item(27, 3).Value = 'hc_sochi'
$ws1.Cells.Item(27, 4).Value = 'Хёфенмайер Ноэль'
$ws1.Cells.Item(27, 5).Value = '22'
$ws1.Cells.Item(27, 6).Value = 'защитник'
$ws1.Cells.Item(27, 7).Value = '44847'
$ws1.Cells.Item(27, 8).Value = '1369_СОЧ_хефенмайерноэль'
$ws1.Cells.Item(27, 9).Value = 'injured_active'
$ws1.Cells.Item(27, 10).Value = 'https://www.khl.ru/clubs/hc_sochi/team/'
$ws1.Cells.Item(27, 11).Value = '2025-11-29T03:02:21.208237+00:00'

$ws1.Cells.Item(28, 1).Value = 'СПР'
$ws1.Cells.Item(28, 2).Value = 'Спартак'
$ws1.Cells.Item(28, 3).Value = 'spartak'
$ws1.Cells.Item(28, 4).Value = 'Вишневский Дмитрий'
$ws1.Cells.Item(28, 5).Value = '55'
$ws1.Cells.Item(28, 6).Value = 'защитник'
$ws1.Cells.Item(28, 7).Value = '15299'
$ws1.Cells.Item(28, 8).Value = '1369_СПР_вишневскийдмитрий'
$ws1.Cells.Item(28, 9).Value = 'injured_active'
$ws1.Cells.Item(28, 10).Value = 'https://www.khl.ru/clubs/spartak/team/'
$ws1.Cells.Item(28, 11).Value = '2025-11-29T03:02:23.574665+00:00'

$ws1.Cells.Item(29, 1).Value = 'ТОР'
$ws1.Cells.Item(29, 2).Value = 'Торпедо'
$ws1.Cells.Item(29, 3).Value = 'torpedo'
$ws1.Cells.Item(29, 4).Value = 'Науменков Михаил'
$ws1.Cells.Item(29, 5).Value = '33'
$ws1.Cells.Item(29, 6).Value = 'защитник'
$ws1.Cells.Item(29, 7).Value = '16400'
$ws1.Cells.Item(29, 8).Value = '1369_ТОР_науменковмихаил'
$ws1.Cells.Item(29, 9).Value = 'injured_active'
$ws1.Cells.Item(29, 10).Value = 'https://www.khl.ru/clubs/torpedo/team/'
$ws1.Cells.Item(29, 11).Value = '2025-11-29T03:02:28.584902+00:00'

$ws1.Cells.Item(30, 1).Value = 'ЦСК'
$ws1.Cells.Item(30, 2).Value = 'ЦСКА'
$ws1.Cells.Item(30, 3).Value = 'cska'
$ws1.Cells.Item(30, 4).Value = 'Бучельников Дмитрий'
$ws1.Cells.Item(30, 5).Value = '72'
$ws1.Cells.Item(30, 6).Value = 'нападающий'
$ws1.Cells.Item(30, 7).Value = '39102'
$ws1.Cells.Item(30, 8).Value = '1369_ЦСК_бучельниковдмитрий'
$ws1.Cells.Item(30, 9).Value = 'injured_active'
$ws1.Cells.Item(30, 10).Value = 'https://www.khl.ru/clubs/cska/team/'
$ws1.Cells.Item(30, 11).Value = '2025-11-29T03:02:34.114741+00:00'

$ws1.Cells.Item(31, 1).Value = 'ЦСК'
$ws1.Cells.Item(31, 2).Value = 'ЦСКА'
$ws1.Cells.Item(31, 3).Value = 'cska'
$ws1.Cells.Item(31, 4).Value = 'Моисеев Данила'
$ws1.Cells.Item(31, 5).Value = '93'
$ws1.Cells.Item(31, 6).Value = 'нападающий'
$ws1.Cells.Item(31, 7).Value = '23931'
$ws1.Cells.Item(31, 8).Value = '1369_ЦСК_моисеевданила'
$ws1.Cells.Item(31, 9).Value = 'injured_active'
$ws1.Cells.Item(31, 10).Value = 'https://www.khl.ru/clubs/cska/team/'
$ws1.Cells.Item(31, 11).Value = '2025-11-29T03:02:34.114772+00:00'

$ws1.Cells.Item(32, 1).Value = 'ШДР'
$ws1.Cells.Item(32, 2).Value = 'Драконы'
$ws1.Cells.Item(32, 3).Value = 'kunlun'
$ws1.Cells.Item(32, 4).Value = 'Бишофф Джейк'
$ws1.Cells.Item(32, 5).Value = '28'
$ws1.Cells.Item(32, 6).Value = 'защитник'
$ws1.Cells.Item(32, 7).Value = '45490'
$ws1.Cells.Item(32, 8).Value = '1369_ШДР_бишоффджейк'
$ws1.Cells.Item(32, 9).Value = 'injured_active'
$ws1.Cells.Item(32, 10).Value = 'https://www.khl.ru/clubs/kunlun/team/'
$ws1.Cells.Item(32, 11).Value = '2025-11-29T03:02:36.829095+00:00'

$ws1.Cells.Item(33, 1).Value = 'ШДР'
$ws1.Cells.Item(33, 2).Value = 'Драконы'
$ws1.Cells.Item(33, 3).Value = 'kunlun'
$ws1.Cells.Item(33, 4).Value = 'Гроло Жереми'
$ws1.Cells.Item(33, 5).Value = '75'
$ws1.Cells.Item(33, 6).Value = 'защитник'
$ws1.Cells.Item(33, 7).Value = '45343'
$ws1.Cells.Item(33, 8).Value = '1369_ШДР_гроложереми'
$ws1.Cells.Item(33, 9).Value = 'injured_active'
$ws1.Cells.Item(33, 10).Value = 'https://www.khl.ru/clubs/kunlun/team/'
$ws1.Cells.Item(33, 11).Value = '2025-11-29T03:02:36.829129+00:00'

# --- Clear old row 34 on sheet1 (previously had data, now removed) ---
$ws1.Cells.Item(34, 1).Value = ""
$ws1.Cells.Item(34, 2).Value = ""
$ws1.Cells.Item(34, 3).Value = ""
$ws1.Cells.Item(34, 4).Value = ""
$ws1.Cells.Item(34, 5).Value = ""
$ws1.Cells.Item(34, 6).Value = ""
$ws1.Cells.Item(34, 7).Value = ""
$ws1.Cells.Item(34, 8).Value = ""
$ws1.Cells.Item(34, 9).Value = ""
$ws1.Cells.Item(34, 10).Value = ""
$ws1.Cells.Item(34, 11).Value = ""

# --- Sheet2 "returned": rewrite rows 2..4 ---
$ws2.Cells.Item(2, 1).Value = 'СИБ'
$ws2.Cells.Item(2, 2).Value = 'Сибирь'
$ws2.Cells.Item(2, 3).Value = 'Приски Чейз Эванс'
$ws2.Cells.Item(2, 4).Value = '1369_СИБ_прискичейзэванс'
$ws2.Cells.Item(2, 5).Value = 'RETURN'
$ws2.Cells.Item(2, 6).Value = '2025-11-29T11:02:37.340858+08:00'
$ws2.Cells.Item(2, 7).Value = '2025-11-29'

$ws2.Cells.Item(3, 1).Value = 'СОЧ'
$ws2.Cells.Item(3, 2).Value = 'ХК Сочи'
$ws2.Cells.Item(3, 3).Value = 'Гуськов Матвей'
$ws2.Cells.Item(3, 4).Value = '1369_СОЧ_гуськовматвей'
$ws2.Cells.Item(3, 5).Value = 'RETURN'
$ws2.Cells.Item(3, 6).Value = '2025-11-29T11:02:37.340858+08:00'
$ws2.Cells.Item(3, 7).Value = '2025-11-29'

$ws2.Cells.Item(4, 1).Value = 'СЮЛ'
$ws2.Cells.Item(4, 2).Value = 'Салават Юлаев'
$ws2.Cells.Item(4, 3).Value = 'Хворов Николай'
$ws2.Cells.Item(4, 4).Value = '1369_СЮЛ_хворовниколай'
$ws2.Cells.Item(4, 5).Value = 'RETURN'
$ws2.Cells.Item(4, 6).Value = '2025-11-29T11:02:37.340858+08:00'
$ws2.Cells.Item(4, 7).Value = '2025-11-29'

# --- Sheet3 "new_injured": add rows 2..3 ---
$ws3.Cells.Item(2, 1).Value = 'ММГ'
$ws3.Cells.Item(2, 2).Value = 'Металлург Мг'
$ws3.Cells.Item(2, 3).Value = 'Сиряцкий Александр'
$ws3.Cells.Item(2, 4).Value = '1369_ММГ_сиряцкийалександр'
$ws3.Cells.Item(2, 5).Value = 'INJURED_NEW'
$ws3.Cells.Item(2, 6).Value = '2025-11-29T11:02:37.340858+08:00'
$ws3.Cells.Item(2, 7).Value = '2025-11-29'

$ws3.Cells.Item(3, 1).Value = 'СОЧ'
$ws3.Cells.Item(3, 2).Value = 'ХК Сочи'
$ws3.Cells.Item(3, 3).Value = 'Бикмуллин Рафаэль'
$ws3.Cells.Item(3, 4).Value = '1369_СОЧ_бикмуллинрафаэль'
$ws3.Cells.Item(3, 5).Value = 'INJURED_NEW'
$ws3.Cells.Item(3, 6).Value = '2025-11-29T11:02:37.340858+08:00'
$ws3.Cells.Item(3, 7).Value = '2025-11-29'

